$d = $word.ActiveDocument

# Remove the first nine paragraphs ("1", "22", "333", "4444",
# "55555 commit", "713713713713", "Master add", the blank paragraph,
# and "Add conflicts content"), leaving the "714delete" paragraph as
# the new first paragraph.
$r = $d.Range(0, $d.Paragraphs.Item(10).Range.Start)
$r.Delete()

# Turn "714delete" into "Master had a bug" in place, preserving the
# run formatting and the _GoBack bookmark that follows it.
$d.Content.Find.Execute("714delete", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Master had a bug", 2)

# Collapse the two trailing empty paragraphs into one.
$d.Paragraphs.Item(2).Range.Delete()
